$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 126
$ws.Range("J2").Value = 447
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 143
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 84
$ws.Range("O2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 56
$ws.Range("T2").Value = 99
$ws.Range("V2").Value = 706
$ws.Range("X2").Value = 692
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 7
$ws.Range("AA2").Value = 5

$wb.Save()
